$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----
$ws.Range("A1").Value = "uID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "sID"
$ws.Range("E1").Value = "Notes"

# ---- Shared formatting for the uID column (right aligned, Arial 10) ----
$ws.Range("A2:A4").Font.Name = "Arial"
$ws.Range("A2:A4").Font.Size = 10
$ws.Range("A2:A4").HorizontalAlignment = -4152

# ---- Shared formatting for the Name/sID columns (left aligned, Arial 10) ----
$ws.Range("B2:B4").Font.Name = "Arial"
$ws.Range("B2:B4").Font.Size = 10
$ws.Range("B2:B4").HorizontalAlignment = -4131

$ws.Range("D2:D4").Font.Name = "Arial"
$ws.Range("D2:D4").Font.Size = 10
$ws.Range("D2:D4").HorizontalAlignment = -4131

# ---- Date column kept as literal text "2025-12-11" (force text so it is
#      not auto-converted into a date serial number) ----
$ws.Range("C2:C4").NumberFormat = "@"
$ws.Range("C3:C4").Font.Name = "Arial"
$ws.Range("C3:C4").Font.Size = 10
$ws.Range("C3:C4").HorizontalAlignment = -4152

# Row 2
$ws.Range("A2").Value = 66001
$ws.Range("B2").Value = "A"
$ws.Range("C2").Value = "2025-12-11"
$ws.Range("D2").Value = "S1"

# Row 3
$ws.Range("A3").Value = 66001
$ws.Range("B3").Value = "A"
$ws.Range("C3").Value = "2025-12-11"
$ws.Range("D3").Value = "S2"

# Row 4
$ws.Range("A4").Value = 66002
$ws.Range("B4").Value = "B"
$ws.Range("C4").Value = "2025-12-11"
$ws.Range("D4").Value = "S1"

# ---- Row heights (matching data rows) ----
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 15

# ---- View state: leave selection on G13, as in the source workbook ----
$ws.Range("G13").Select()
